$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price / 1h-volume figures.
# Each target cell is forced to Text storage (NumberFormat "@") before the
# write so Excel does not auto-convert number-looking strings (e.g. "27.65")
# into numeric cells, then the style is reset back to "Normal" so no stray
# cell-format (s="...") is left behind relative to the original workbook.
$updates = @(
    @{ Addr = "D2"; Value = "63.651.46" },
    @{ Addr = "E2"; Value = "  +0.17%  " },
    @{ Addr = "D3"; Value = "2.651.04" },
    @{ Addr = "E3"; Value = "  +0.00%  " },
    @{ Addr = "E4"; Value = "  -0.02%  " },
    @{ Addr = "D5"; Value = "603.68" },
    @{ Addr = "E5"; Value = "  +2.02%  " },
    @{ Addr = "D6"; Value = "147.03" },
    @{ Addr = "E6"; Value = "  +2.00%  " },
    @{ Addr = "E8"; Value = "  +0.31%  " },
    @{ Addr = "E9"; Value = "  +1.59%  " },
    @{ Addr = "E10"; Value = "  -0.43%  " },
    @{ Addr = "E11"; Value = "  +4.63%  " },
    @{ Addr = "E12"; Value = "  -0.14%  " },
    @{ Addr = "D13"; Value = "27.65" },
    @{ Addr = "E13"; Value = "  +0.54%  " },
    @{ Addr = "D14"; Value = "3.125.49" },
    @{ Addr = "E14"; Value = "  +0.01%  " },
    @{ Addr = "D15"; Value = "63.495.47" },
    @{ Addr = "E15"; Value = "  +0.08%  " },
    @{ Addr = "E16"; Value = "  +0.99%  " },
    @{ Addr = "D17"; Value = "2.653.33" },
    @{ Addr = "E17"; Value = "  +0.36%  " },
    @{ Addr = "D18"; Value = "11.54" },
    @{ Addr = "E18"; Value = "  +0.93%  " },
    @{ Addr = "D19"; Value = "4.57" },
    @{ Addr = "E19"; Value = "  +4.60%  " },
    @{ Addr = "D20"; Value = "343.51" },
    @{ Addr = "E20"; Value = "  +0.87%  " },
    @{ Addr = "D21"; Value = "6.95" },
    @{ Addr = "E21"; Value = "  +3.09%  " },
    @{ Addr = "E22"; Value = "  -0.08%  " },
    @{ Addr = "E23"; Value = "  -3.02%  " },
    @{ Addr = "D24"; Value = "66.93" },
    @{ Addr = "E24"; Value = "  -0.38%  " },
    @{ Addr = "D25"; Value = "1.71" },
    @{ Addr = "E25"; Value = "  +1.85%  " },
    @{ Addr = "D26"; Value = "9.07" },
    @{ Addr = "E26"; Value = "  +7.62%  " },
    @{ Addr = "D27"; Value = "574.83" },
    @{ Addr = "E27"; Value = "  +5.81%  " },
    @{ Addr = "D28"; Value = "1.57" },
    @{ Addr = "E28"; Value = "  +1.31%  " },
    @{ Addr = "E29"; Value = "  -1.67%  " },
    @{ Addr = "D30"; Value = "8.01" },
    @{ Addr = "E30"; Value = "  +2.86%  " },
    @{ Addr = "E31"; Value = "  -0.10%  " },
    @{ Addr = "E32"; Value = "  +3.77%  " },
    @{ Addr = "E33"; Value = "  -3.56%  " },
    @{ Addr = "D34"; Value = "0.0₃0825" },
    @{ Addr = "E34"; Value = "  +2.09%  " },
    @{ Addr = "D35"; Value = "5.23" },
    @{ Addr = "E35"; Value = "  +6.86%  " },
    @{ Addr = "D36"; Value = "168.80" },
    @{ Addr = "E36"; Value = "  -3.77%  " },
    @{ Addr = "E37"; Value = "  +1.21%  " },
    @{ Addr = "E38"; Value = "  -0.06%  " },
    @{ Addr = "E39"; Value = "  +7.09%  " },
    @{ Addr = "D40"; Value = "19.13" },
    @{ Addr = "E40"; Value = "  +0.19%  " },
    @{ Addr = "E41"; Value = "  +0.01%  " },
    @{ Addr = "D42"; Value = "169.26" },
    @{ Addr = "E42"; Value = "  -0.93%  " },
    @{ Addr = "E43"; Value = "  +1.03%  " },
    @{ Addr = "D44"; Value = "22.23" },
    @{ Addr = "E44"; Value = "  -1.21%  " },
    @{ Addr = "D45"; Value = "0.0573" },
    @{ Addr = "E45"; Value = "  +2.57%  " },
    @{ Addr = "E46"; Value = "  +0.18%  " },
    @{ Addr = "E47"; Value = "  +3.09%  " },
    @{ Addr = "D48"; Value = "0.0963" },
    @{ Addr = "E48"; Value = "  +0.17%  " },
    @{ Addr = "D49"; Value = "1.90" },
    @{ Addr = "E49"; Value = "  +10.89%  " },
    @{ Addr = "D50"; Value = "18.87" },
    @{ Addr = "E50"; Value = "  +0.11%  " },
    @{ Addr = "D51"; Value = "0.178" },
    @{ Addr = "E51"; Value = "  +1.84%  " }
)

foreach ($item in $updates) {
    $cell = $ws.Range($item.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}
